$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the old row 24 ("Flame Seedless",
# Provincia del Elquí, 2021-12-16), pushing all rows from 24 down to 25-58.
$ws.Rows.Item(24).Insert()

$ws.Range("A24").Value = 8
$ws.Range("B24").Value = "Terminal La Palmera de La Serena"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44546
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100109
$ws.Range("H24").Value = "Uva"
$ws.Range("I24").Value = 100109001
$ws.Range("J24").Value = "Uva"
$ws.Range("K24").Value = "Flame Seedless"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 400
$ws.Range("N24").Value = 11500
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 11750
$ws.Range("Q24").Value = "`$/bandeja 18 kilos"
$ws.Range("R24").Value = "Provincia del Elquí"
$ws.Range("S24").Value = 653
$ws.Range("T24").Value = 18
